$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 443; this shifts the existing rows
# 443..526 down to 444..527 and grows the sheet dimension to R527.
$ws.Rows.Item(443).Insert()

# Populate the newly inserted row 443 with the new weekly record.
$ws.Range("A443").Value = 3
$ws.Range("B443").Value = "Femacal de La Calera"
$ws.Range("C443").Value = "Coquimbo"
$ws.Range("D443").Value = 44995
$ws.Range("E443").Value = 5
$ws.Range("F443").Value = 100112040
$ws.Range("G443").Value = "Cilantro"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 145
$ws.Range("K443").Value = 6000
$ws.Range("L443").Value = 6500
$ws.Range("M443").Value = 6259
$ws.Range("N443").Value = '$/docena de atados (3 kilos)'
$ws.Range("O443").Value = "Provincia de Quillota"
$ws.Range("P443").Value = 2086
$ws.Range("Q443").Value = 3
$ws.Range("R443").Value = "Hortaliza"
